# Coding cleaning script to fix up the raw data from the experts.
# Replaces the raw numeric adcap_score values in column D with categorical
# labels (none/low/medium/high) on a few of the trait_stressor sheets, adds
# a missing note, and leaves "poisons + toxins" as the active/selected sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: biomass removal -> D2:D15 all become "none"
# ---------------------------------------------------------------------
$wsBiomass = $wb.Worksheets.Item("biomass removal")
$wsBiomass.Activate()
$wsBiomass.Range("D2").Value = "none"
$wsBiomass.Range("D3").Value = "none"
$wsBiomass.Range("D4").Value = "none"
$wsBiomass.Range("D5").Value = "none"
$wsBiomass.Range("D6").Value = "none"
$wsBiomass.Range("D7").Value = "none"
$wsBiomass.Range("D8").Value = "none"
$wsBiomass.Range("D9").Value = "none"
$wsBiomass.Range("D10").Value = "none"
$wsBiomass.Range("D11").Value = "none"
$wsBiomass.Range("D12").Value = "none"
$wsBiomass.Range("D13").Value = "none"
$wsBiomass.Range("D14").Value = "none"
$wsBiomass.Range("D15").Value = "none"
$wsBiomass.Range("D2:D8").Select()

# ---------------------------------------------------------------------
# Sheet: habitat loss + degradation
#   D2:D3   -> "none"
#   D4      -> "low"
#   D5      -> "medium"
#   D6:D8   -> "high"
#   D10:D15 -> 0 (stay numeric)
# ---------------------------------------------------------------------
$wsHabitat = $wb.Worksheets.Item("habitat loss + degradation")
$wsHabitat.Activate()
$wsHabitat.Range("D2").Value = "none"
$wsHabitat.Range("D3").Value = "none"
$wsHabitat.Range("D4").Value = "low"
$wsHabitat.Range("D5").Value = "medium"
$wsHabitat.Range("D6").Value = "high"
$wsHabitat.Range("D7").Value = "high"
$wsHabitat.Range("D8").Value = "high"
$wsHabitat.Range("D10").Value = 0
$wsHabitat.Range("D11").Value = 0
$wsHabitat.Range("D12").Value = 0
$wsHabitat.Range("D13").Value = 0
$wsHabitat.Range("D14").Value = 0
$wsHabitat.Range("D15").Value = 0
$wsHabitat.Range("D9:D15").Select()

# ---------------------------------------------------------------------
# Sheet: poisons + toxins
#   D2      -> "none"
#   D9      -> "high"   (+ note added in E9)
#   D10:D11 -> "medium"
#   D12:D13 -> "low"
#   D14     -> "none"
#   D15     -> "high"
# ---------------------------------------------------------------------
$wsPoisons = $wb.Worksheets.Item("poisons + toxins")
$wsPoisons.Activate()
$wsPoisons.Range("D2").Value = "none"
$wsPoisons.Range("D9").Value = "high"
$wsPoisons.Range("E9").Value = "assume short PLD means less exposure to poisons and toxins"
$wsPoisons.Range("D10").Value = "medium"
$wsPoisons.Range("D11").Value = "medium"
$wsPoisons.Range("D12").Value = "low"
$wsPoisons.Range("D13").Value = "low"
$wsPoisons.Range("D14").Value = "none"
$wsPoisons.Range("D15").Value = "high"
$wsPoisons.Range("D3").Select()
